$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 373, pushing the existing row 373 (and everything
# below it) down by one. This extends the used range from R485 to R486.
$ws.Rows.Item(373).Insert()

# Populate the newly inserted row 373 with the new weekly price record.
$ws.Cells.Item(373, 1).Value = 5
$ws.Cells.Item(373, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(373, 3).Value = "Maule"
$ws.Cells.Item(373, 4).Value = 44985
$ws.Cells.Item(373, 5).Value = 7
$ws.Cells.Item(373, 6).Value = 100114014
$ws.Cells.Item(373, 7).Value = "Betarraga"
$ws.Cells.Item(373, 8).Value = "Sin especificar"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 6000
$ws.Cells.Item(373, 11).Value = 500
$ws.Cells.Item(373, 12).Value = 550
$ws.Cells.Item(373, 13).Value = 525
$ws.Cells.Item(373, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(373, 15).Value = "Región del Maule"
$ws.Cells.Item(373, 16).Value = 105
$ws.Cells.Item(373, 17).Value = 5
$ws.Cells.Item(373, 18).Value = "Hortaliza"
